$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 30.02419053477981
$ws.Range("D2").Value = 0.1041905347798036
$ws.Range("E2").Value = 0.01085566753770145
$ws.Range("C3").Value = 30.02977761739633
$ws.Range("D3").Value = 0.04977761739632669
$ws.Range("E3").Value = 0.002477811193655085
$ws.Range("C4").Value = 30.07035366584773
$ws.Range("D4").Value = 0.03035366584772703
$ws.Range("E4").Value = 0.0009213450303954703
$ws.Range("C5").Value = 30.08315331313932
$ws.Range("D5").Value = -0.1268466868606843
$ws.Range("E5").Value = 0.01609008196753251
$ws.Range("C6").Value = 30.19825301488454
$ws.Range("D6").Value = -0.02174698511546325
$ws.Range("E6").Value = 0.00047293136161218
$ws.Range("C7").Value = 30.25568284039331
$ws.Range("D7").Value = -0.1243171596066937
$ws.Range("E7").Value = 0.01545475617267614
$ws.Range("C8").Value = 30.54227831882804
$ws.Range("D8").Value = 0.1022783188280414
$ws.Range("E8").Value = 0.01046085450229049
$ws.Range("C9").Value = 30.5519949621441
$ws.Range("D9").Value = 0.07199496214410317
$ws.Range("E9").Value = 0.005183274574130849
$ws.Range("C10").Value = 30.60082321928214
$ws.Range("D10").Value = -0.08917678071786383
$ws.Range("E10").Value = 0.007952498219201971
$ws.Range("C11").Value = 30.76713189625077
$ws.Range("D11").Value = 0.01713189625076694
$ws.Range("E11").Value = 0.0002935018691470424
$ws.Range("C12").Value = 30.89341029886786
$ws.Range("D12").Value = -0.04658970113214522
$ws.Range("E12").Value = 0.002170600251582614
$ws.Range("C13").Value = 31.10379906852334
$ws.Range("D13").Value = 0.1537990685233375
$ws.Range("E13").Value = 0.02365415347864626
$ws.Range("C14").Value = 31.27142860609367
$ws.Range("D14").Value = 0.2514286060936683
$ws.Range("E14").Value = 0.06321634396220505
$ws.Range("C15").Value = 31.33366433097254
$ws.Range("D15").Value = 0.2136643309725379
$ws.Range("E15").Value = 0.04565244632994222
$ws.Range("C16").Value = 31.1674750055018
$ws.Range("D16").Value = -0.1125249944982052
$ws.Range("E16").Value = 0.01266187438682111
$ws.Range("C17").Value = 30.93389696060447
$ws.Range("D17").Value = -0.4461030393955276
$ws.Range("E17").Value = 0.1990079217579277
$ws.Range("C18").Value = 31.12874365607843
$ws.Range("D18").Value = -0.4512563439215675
$ws.Range("E18").Value = 0.2036322879294601
$ws.Range("C19").Value = 31.68285551433151
$ws.Range("D19").Value = 0.03285551433151213
$ws.Range("E19").Value = 0.001079484821988199
$ws.Range("C20").Value = 32.54857965831769
$ws.Range("D20").Value = 0.6685796583176931
$ws.Range("E20").Value = 0.4469987595162032
$ws.Range("C21").Value = 32.55971918219023
$ws.Range("D21").Value = 0.2797191821902274
$ws.Range("E21").Value = 0.07824282088516962
$ws.Range("C22").Value = 32.73532906098168
$ws.Range("D22").Value = 0.2853290609816739
$ws.Range("E22").Value = 0.08141267304068377
$ws.Range("C23").Value = 32.9375852778052
$ws.Range("D23").Value = 0.08758527780519643
$ws.Range("E23").Value = 0.007671180888213434
$ws.Range("C24").Value = 32.99351977061944
$ws.Range("D24").Value = 0.09351977061943728
$ws.Range("E24").Value = 0.008745947496712165
$ws.Range("C25").Value = 33.22201537420774
$ws.Range("D25").Value = 0.1220153742077414
$ws.Range("E25").Value = 0.01488775154305518
$ws.Range("C26").Value = 33.51113863081579
$ws.Range("D26").Value = 0.1111386308157876
$ws.Range("E26").Value = 0.01235179525960793
$ws.Range("C27").Value = 33.60283648162176
$ws.Range("D27").Value = -0.09716351837824533
$ws.Range("E27").Value = 0.009440749303639617
$ws.Range("C28").Value = 33.71172498514026
$ws.Range("D28").Value = -0.3882750148597438
$ws.Range("E28").Value = 0.1507574871643343
$ws.Range("C29").Value = 34.29938866808003
$ws.Range("D29").Value = -0.1006113319199713
$ws.Range("E29").Value = 0.01012264011071063
$ws.Range("C30").Value = 34.59425327533567
$ws.Range("D30").Value = -0.3057467246643242
$ws.Range("E30").Value = 0.09348105964296206
$ws.Range("C31").Value = 35.63727115299272
$ws.Range("D31").Value = 0.3372711529927201
$ws.Range("E31").Value = 0.1137518306410388
$ws.Range("C32").Value = 35.89773605532633
$ws.Range("D32").Value = 0.1977360553263239
$ws.Range("E32").Value = 0.03909954757601502
$ws.Range("C33").Value = 36.32664756839336
$ws.Range("D33").Value = 0.02664756839336491
$ws.Range("E33").Value = 0.0007100929012790603
$ws.Range("C34").Value = 36.82917752352359
$ws.Range("D34").Value = 0.02917752352359315
$ws.Range("E34").Value = 0.0008513278789698317
$ws.Range("C35").Value = 36.99573242149806
$ws.Range("D35").Value = -0.3042675785019355
$ws.Range("E35").Value = 0.0925787593274315
$ws.Range("C36").Value = 37.78072258726555
$ws.Range("D36").Value = -0.1192774127344478
$ws.Range("E36").Value = 0.01422710118862381
$ws.Range("C37").Value = 38.52270208870895
$ws.Range("D37").Value = 0.02270208870894663
$ws.Range("E37").Value = 0.0005153848317488821
$ws.Range("C38").Value = 39.12314879822863
$ws.Range("D38").Value = 0.2231487982286282
$ws.Range("E38").Value = 0.04979538615088103
$ws.Range("C39").Value = 39.44146703013038
$ws.Range("D39").Value = 0.04146703013037722
$ws.Range("E39").Value = 0.001719514587833613
$ws.Range("C40").Value = 39.72263802358091
$ws.Range("D40").Value = -0.1773619764190926
$ws.Range("E40").Value = 0.03145727067928675
$ws.Range("C41").Value = 39.7702638787822
$ws.Range("D41").Value = -0.3297361212177989
$ws.Range("E41").Value = 0.108725909635759
$ws.Range("C42").Value = 39.90697336106128
$ws.Range("D42").Value = -0.6930266389387256
$ws.Range("E42").Value = 0.4802859222787068
$ws.Range("C43").Value = 40.19635090321012
$ws.Range("D43").Value = -0.7036490967898743
$ws.Range("E43").Value = 0.4951220514132059
$ws.Range("C44").Value = 41.21758245333769
$ws.Range("D44").Value = 0.01758245333768826
$ws.Range("E44").Value = 0.0003091426653719852
$ws.Range("C45").Value = 41.28955640879663
$ws.Range("D45").Value = -0.2104435912033651
$ws.Range("E45").Value = 0.04428650507856904
$ws.Range("C46").Value = 41.71509244587822
$ws.Range("D46").Value = -0.0849075541217772
$ws.Range("E46").Value = 0.007209292746942523
$ws.Range("C47").Value = 42.40166043411897
$ws.Range("D47").Value = 0.2016604341189705
$ws.Range("E47").Value = 0.04066693068905164
$ws.Range("C48").Value = 43.53255751544864
$ws.Range("D48").Value = 0.8325575154486344
$ws.Range("E48").Value = 0.6931520165300031
$ws.Range("C49").Value = 43.95549435339204
$ws.Range("D49").Value = 0.255494353392038
$ws.Range("E49").Value = 0.06527736461521558
$ws.Range("C50").Value = 44.25777747149995
$ws.Range("D50").Value = 0.05777747149994639
$ws.Range("E50").Value = 0.003338236212927117
$ws.Range("C51").Value = 45.62681281157661
$ws.Range("D51").Value = 0.02681281157660464
$ws.Range("E51").Value = 0.000718926864642504
$ws.Range("C52").Value = 0.01236847578596567
$ws.Range("E52").Value = 3.819149214691712
$ws.Range("E53").Value = 0.07638298429383424
